# Update cryptos list with latest scraped prices/volumes (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings like "26.381.41" that Excel would otherwise
# auto-convert to numbers (losing the thousands-dot formatting / precision).
# Force the whole price column to Text before writing, then restore the
# default "Normal" style so no stray number-format sticks to the cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.381.41"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.593.21"
$ws.Range("E3").Value = "  +0.55%  "

# Row 5 - BNB
$ws.Range("D5").Value = "211.64"
$ws.Range("E5").Value = "  +0.89%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.19%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.36%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.53%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.20%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -0.52%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.21%  "

# Row 12 - Wrapped liquid staked Ether 2.0
$ws.Range("E12").Value = "  +0.50%  "

# Row 13 - Wrapped Ether
$ws.Range("D13").Value = "1.606.77"
$ws.Range("E13").Value = "  +1.48%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "4.07"
$ws.Range("E14").Value = "  +1.30%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.524"

# Row 16 - Litecoin
$ws.Range("D16").Value = "64.79"
$ws.Range("E16").Value = "  +0.33%  "

# Row 17 - Wrapped BTC
$ws.Range("D17").Value = "26.370.58"

# Row 18 - Shiba Inu
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -0.98%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "7.52"
$ws.Range("E19").Value = "  +4.35%  "

# Row 20 - Bitcoin Cash
$ws.Range("D20").Value = "212.17"
$ws.Range("E20").Value = "  +2.55%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.33%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.29%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.04"
$ws.Range("E23").Value = "  +2.19%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.77%  "

# Row 25 - Monero
$ws.Range("D25").Value = "143.94"

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.25%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "7.09"
$ws.Range("E27").Value = "  +0.96%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.74%  "

# Row 29 - Ethereum Classic
$ws.Range("E29").Value = "  -0.03%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.11%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.03%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.34%  "

# Row 33 - Internet Computer (DFINITY)
$ws.Range("E33").Value = "  +1.29%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.338.33"
$ws.Range("E34").Value = "  +4.17%  "

# Row 35 - Huobi Token
$ws.Range("E35").Value = "  -1.29%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -0.85%  "

# Row 37 - Lido DAO Token
$ws.Range("E37").Value = "  +0.04%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.15%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +0.20%  "

# Row 40 - Frax Share
$ws.Range("D40").Value = "5.78"
$ws.Range("E40").Value = "  +5.71%  "

# Row 41 - Pax Dollar
$ws.Range("E41").Value = "  -0.35%  "

# Row 42 - WEMIX Token
$ws.Range("D42").Value = "1.02"
$ws.Range("E42").Value = "  -22.16%  "

# Rows 43/44 swap: TrustWalletToken now ranks above MXToken
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.767"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45 - Rocket Pool ETH
$ws.Range("D45").Value = "1.729.76"
$ws.Range("E45").Value = "  +0.54%  "

# Row 46 - Aave
$ws.Range("E46").Value = "  -0.59%  "

# Row 47 - Quant
$ws.Range("D47").Value = "88.13"

# Row 48 - Render Token
$ws.Range("E48").Value = "  -3.71%  "

# Row 49 - Algorand
$ws.Range("D49").Value = "0.0984"
$ws.Range("E49").Value = "  -3.60%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -0.91%  "

# Row 51 - USDD
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.45%  "

# Restore the default style on the price column so no cell keeps a stray
# explicit text-number-format that wasn't there before.
$ws.Range("D2:D51").Style = "Normal"
